$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "MAE" column is being inserted before the existing "Tipo" column.
# First, move the "Tipo" header (together with its header formatting) from
# D1 over to the new E1 position.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").Value2 = $ws.Range("D1").Value2

# Move the corresponding data value ("single") from D2 to E2.
$ws.Range("E2").Value2 = $ws.Range("D2").Value2

# Put the new "MAE" header into D1 (it keeps the header style already there).
$ws.Range("D1").Value2 = "MAE"

# Refresh the computed metrics (MSE, R2) and add the new MAE value.
$ws.Range("B2").Value2 = 0.4919149660744245
$ws.Range("C2").Value2 = 0.9855511688435743
$ws.Range("D2").Value2 = 0.5134388839929866
